# Fixes #1046 use 2 template dictionaries with and without LLOQ (#1273)
#
# The "tpDictionary" sheet is duplicated into a new sheet "tpDictionaryLoq"
# placed right after it. The new sheet keeps the LLOQ row (row 12) while the
# original "tpDictionary" sheet drops that row (it now only documents
# parameters that are used without a LLOQ column).

$wb = $excel.ActiveWorkbook

$original = $wb.Worksheets.Item("tpDictionary")

# Duplicate "tpDictionary" right after itself.
$original.Copy([System.Reflection.Missing]::Value, $original)
$loq = $wb.Worksheets.Item("tpDictionary (2)")
$loq.Name = "tpDictionaryLoq"

# The duplicated sheet's cells inherited the wrap-text alignment from the
# source sheet; the published version turns that off.
$loq.Range("A1:G12").WrapText = $false

# The original sheet no longer needs the LLOQ-specific row (former row 12:
# "67"/"60"/"68"/.../"64").
$original.Rows(12).Delete()

# Restore a sensible selection/active-sheet state on both sheets.
$original.Range("C17").Select() | Out-Null
$loq.Range("C15").Select() | Out-Null
$loq.Activate() | Out-Null
